$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.898.36"
$ws.Range("E2").Value = "  -0.48%  "

$ws.Range("D3").Value = "'2.300.11"
$ws.Range("E3").Value = "  -0.29%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'305.29"
$ws.Range("E5").Value = "  +1.52%  "

$ws.Range("D6").Value = "'97.35"
$ws.Range("E6").Value = "  -0.52%  "

$ws.Range("E7").Value = "  -1.57%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("E9").Value = "  -2.55%  "

$ws.Range("E10").Value = "  -0.47%  "

$ws.Range("E11").Value = "  +0.70%  "

$ws.Range("D12").Value = "'18.32"
$ws.Range("E12").Value = "  +1.97%  "

$ws.Range("E13").Value = "  +0.98%  "

$ws.Range("E14").Value = "  -1.37%  "

$ws.Range("D15").Value = "'2.658.66"
$ws.Range("E15").Value = "  -0.23%  "

$ws.Range("D16").Value = "'2.300.51"
$ws.Range("E16").Value = "  +1.10%  "

$ws.Range("E17").Value = "  -0.78%  "

$ws.Range("D18").Value = "'42.841.31"
$ws.Range("E18").Value = "  -0.40%  "

$ws.Range("D19").Value = "'13.04"
$ws.Range("E19").Value = "  -1.10%  "

$ws.Range("D20").Value = "'0.0₃0906"
$ws.Range("E20").Value = "  -0.42%  "

$ws.Range("E21").Value = "  -1.19%  "

$ws.Range("D22").Value = "'67.57"
$ws.Range("E22").Value = "  -1.12%  "

$ws.Range("D23").Value = "'236.60"
$ws.Range("E23").Value = "  -0.66%  "

$ws.Range("D24").Value = "'2.18"
$ws.Range("E24").Value = "  -1.32%  "

$ws.Range("E25").Value = "  +2.44%  "

$ws.Range("E26").Value = "  +1.04%  "

$ws.Range("E27").Value = "  +0.25%  "

$ws.Range("D28").Value = "'25.51"
$ws.Range("E28").Value = "  +1.28%  "

$ws.Range("D29").Value = "'167.07"
$ws.Range("E29").Value = "  +0.11%  "

$ws.Range("E30").Value = "  +1.09%  "

$ws.Range("E31").Value = "  -0.83%  "

$ws.Range("D32").Value = "'32.99"
$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("E33").Value = "  +0.08%  "

$ws.Range("D34").Value = "'4.80"
$ws.Range("E34").Value = "  +1.43%  "

$ws.Range("E35").Value = "  -2.15%  "

$ws.Range("D36").Value = "'17.35"
$ws.Range("E36").Value = "  -5.18%  "

$ws.Range("E37").Value = "  -0.34%  "

$ws.Range("E38").Value = "  +0.47%  "

$ws.Range("E39").Value = "  -0.42%  "

$ws.Range("E40").Value = "  -1.90%  "

$ws.Range("E41").Value = "  -1.22%  "

$ws.Range("D42").Value = "'2.73"
$ws.Range("E42").Value = "  -0.72%  "

$ws.Range("D43").Value = "'2.016.24"
$ws.Range("E43").Value = "  +0.38%  "

$ws.Range("E44").Value = "  -1.99%  "

$ws.Range("D45").Value = "'2.14"
$ws.Range("E45").Value = "  -1.24%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'18.08"
$ws.Range("E46").Value = "  +3.14%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'10.04"
$ws.Range("E47").Value = "  -1.66%  "

$ws.Range("D48").Value = "'2.79"
$ws.Range("E48").Value = "  -1.70%  "

$ws.Range("D49").Value = "'2.90"
$ws.Range("E49").Value = "  +7.19%  "

$ws.Range("D50").Value = "'54.01"
$ws.Range("E50").Value = "  -0.90%  "

$ws.Range("D51").Value = "'2.526.12"
$ws.Range("E51").Value = "  -0.48%  "

